$d = $word.ActiveDocument

# 1. Hyperlink timestamp "04:59" -> "12:29" (paragraph 8, first run inside hyperlink)
$d.Content.Find.Execute("04:59", $true, $false, $false, $false, $false, $true, 1, $false, "12:29", 2)

# 2. Highlight body text for Andrés Altamirano -> Mateo De Falco (paragraph 8)
$d.Content.Find.Execute(" ✅: (Andrés Altamirano) Para pensarlo a ver qué es lo del dominio, si le interesa algo más de esto a mí me parecía buena la idea y sobre todo que te dicen analogía que conozco un pibe de Money que nos pide el dueño del money CEO de Money que hizo algo parecido para fondear deportista de polo y como que cuando fondeado los deportistas compraban todo que en todo y después te quedas con un porcentaje de pase, pero los tipos podían llegar a ir a las competencias tener que competir en tal lado y con el fondeo que sería como el mixton, cuál es la campaña que vos hiciste el tipo lograba jugar en un campeonato preparar todas las cosas pero después si hay un mercado de paso el tipo galvania, 2 millones, vos tenés un 10% que lo técnico.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " ✅: (Mateo De Falco) Algo así. Pueden ser dueños de grandes empresas organizaciones y startups?", 2)

# 3. Speaker/time line (paragraph 10)
$d.Content.Find.Execute("00:00 Andrés Altamirano:", $true, $false, $false, $false, $false, $true, 1, $false, "00:00 Mateo De Falco:", 2)

# 4. Dialogue text (paragraph 10)
$d.Content.Find.Execute(" A ver cómo poder hacer todo un circuito en el que me atamos algo de blockchain, o sea de tareas, que si tenemos que hacer esto sería ponerle hoy. Productos para encontrar realmente una aplicación de chelink que vale la pena eso es mínimamente un día hoy por lo menos o mañana que terminar de hoy encontrarlo y terminó de definirlo mañana ya tenés jueves. Después viste que el otro que dijo el", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " innovadores y emprendedores acá innovadores y emprendedores pusimos startupscrito en fase de ideación de desarrollo y de escalamiento plataformas de intercambio nuevas también como exchanges o exchanges o sexys o sea de Nexes o sexys en cuanto a la centralización en este caso Y emprendedores individuales que buscan hacer innovaciones tecnológicas o consultorías independientes es importante diferenciarlos de cualquier tipo de colaborador interno que nosotros tengamos, pero acá también podría llegar a incluirse capaz aquellos emprendedores o innovadores que participan de actividades o dinámicas en el mundo blockchain y podrían verse vinculados sin necesariamente ser colaboradores.", 2)

# 5. Speaker/time line (paragraph 11)
$d.Content.Find.Execute("00:24 Bruno Torossi:", $true, $false, $false, $false, $false, $true, 1, $false, "00:54 Mateo De Falco:", 2)

# 6. Dialogue text (paragraph 11) -- " claro" is ambiguous, scope to this paragraph's Range
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Find.Execute(" claro", $true, $false, $false, $false, $false, $true, 1, $false, " Entonces, cómo cómo piensan? Que podríamos redefinir esta vertical o esta categoría? ideas", 2)

# 7. Speaker/time line (paragraph 12)
$d.Content.Find.Execute("00:25 Andrés Altamirano:", $true, $false, $false, $false, $false, $true, 1, $false, "01:09 Luciano Padovani:", 2)

# 8. Dialogue text (paragraph 12)
$d.Content.Find.Execute(" loco como que también es un producto como andando, aunque sea una transacción, pero ahí me la tiró un poco, qué significa un producto andando en una transacciones?", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " Repetir porque me quedé pensando", 2)

# 9. Speaker/time line (paragraph 13)
$d.Content.Find.Execute("00:33 Gabriel E. Calvo:", $true, $false, $false, $false, $false, $true, 1, $false, "01:11 Mateo De Falco:", 2)

# 10. Dialogue text (paragraph 13)
$d.Content.Find.Execute(" Yo le entendí tratar de es lo que yo decía hay que tratar de completar en una línea finita el flujo completo, o sea, tratar de dar toda la vuelta eso lo que hay que tratar de hacer, digo como que al menos hay que verlo como", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " Ok tenemos acá en la primera", 2)

# 11. Speaker/time line (paragraph 14)
$d.Content.Find.Execute("00:45 Bruno Torossi:", $true, $false, $false, $false, $false, $true, 1, $false, "01:12 Luciano Padovani:", 2)

# 12. Dialogue text (paragraph 14) -- " claro" is ambiguous, scope to this paragraph's Range
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Find.Execute(" claro", $true, $false, $false, $false, $false, $true, 1, $false, " repetir.", 2)

# 13. Speaker/time line (paragraph 15)
$d.Content.Find.Execute("00:46 Gabriel E. Calvo:", $true, $false, $false, $false, $false, $true, 1, $false, "01:13 Mateo De Falco:", 2)

# 14. Dialogue text (paragraph 15)
$d.Content.Find.Execute(" el contrato se ejecuta y termina sucediendo, aunque no tenga un front, aunque vos lo puedas ejecutar el contrato, aunque vos lo pueda otro, lo puede desplollar al contrato, por ejemplo, si tiene front mucho mejor, pero digo eso que fue un poco lo que nosotros en los en los en los bootcam hicimos eso el buscam era copiar el código pegarlo poner el coso conectar la wallete hacer dos boludeces, traer un dato ponerle gas, digo era mucho más simple, no, pero nada que ver, pero nada, entonces eso te permitía bueno ver y ver cómo estaba la transacción registrada, entonces ya.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " categoría innovadores y emprendedores. Para que admito acá tenemos acá innovadores y emprendedores son aquellos actores. Bueno, como lo dice la categoría no sé si hace falta escribirlo dentro de ellos. Yo mencioné startup cripto nuevas plataformas de intercambio y emprendedores individuales. Me gustaría capaz también como no estar capaz contemplado en otras verticales mencionar al famoso criptogró, que participa activamente el ecosistema puede participar en diferentes proyectos blockchain, pero no necesariamente es un colaborador dentro de crecimiento, pero a partir de nuestras activaciones ponerles el día de mañana vamos, si hacemos un evento o colaboramos en un evento definitivamente es under porque es alguien que tiene que participar de las activaciones que hagamos.", 2)

# 15. Speaker/time line (paragraph 16)
$d.Content.Find.Execute("01:28 Gabriel E. Calvo:", $true, $false, $false, $false, $false, $true, 1, $false, "02:02 Luciano Padovani:", 2)

# 16. Dialogue text (paragraph 16)
$d.Content.Find.Execute(" Pero bueno, volviendo al tema que decía lo que como recapitulando. Me parece que es eso digamos de agarrar y palidar con alguien tratar de tener a una reunión como para contarle esto contarle un poco lo que tenemos que problemática es la que vivimos y trata de validarlo con un poco más de ser de certeza, digamos, creo yo que acá hay un foco principal de mi perspectiva que es el foco de de la medición del impacto y nos apareció algo interesante acá, que fue la de poder tomar APIs externa para poder ejecutarlas con un Automation que eso puede estar piola.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " claro, un usuario particular quizás, o sea", 2)

# 17. Speaker/time line (paragraph 17)
$d.Content.Find.Execute("02:01 Gabriel E. Calvo:", $true, $false, $false, $false, $false, $true, 1, $false, "02:06 Mateo De Falco:", 2)

# 18. Dialogue text (paragraph 17)
$d.Content.Find.Execute(" Digamos, por ejemplo consultar a una API de un software, que ya exista, por ejemplo de los que plantan árboles por decirlo o del municipio, por ejemplo del municipio de no sé qué algo, que llega a las plantas o algo así. Conectarlo con algo que tenga que ver con esto y bueno, que uno pueda aportar y que realmente vaya eso sería un golazo bolo, eso poder que vos puedas que en alguna forma podamos hacer que vos le metas plata.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " personas activas en el mundo", 2)

# 19. Speaker/time line (paragraph 18)
$d.Content.Find.Execute("02:32 Gabriel E. Calvo:", $true, $false, $false, $false, $false, $true, 1, $false, "02:16 mariano dueñaz (chat):", 2)

# 20. Dialogue text (paragraph 18)
$d.Content.Find.Execute(" Y que eso automáticamente ejecute acciones en la vida real sería sería un flash, pero bueno, nada de eso", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " power user", 2)

# 21. Speaker/time line (paragraph 19)
$d.Content.Find.Execute("02:41 Andrés Altamirano:", $true, $false, $false, $false, $false, $true, 1, $false, "02:17 Luciano Padovani:", 2)

# 22. Dialogue text (paragraph 19)
$d.Content.Find.Execute(" Se me ocurre como un sitio en", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " sí, sí, pero pero un nombre más", 2)
